$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'auto squeeze',
    'avia active wear women',
    'avia leggings',
    'b and m brown bread',
    'b brown author',
    'b m brown bread',
    'back brace for basketball',
    'back brace for men xl',
    'back brace for workout',
    'back brace full support',
    'back brace hard',
    'back brace running',
    'back brace under clothing',
    'back brace weightlifting',
    'back brace workout',
    'back brace x large',
    'back braces for lower back pain 4x',
    'back compression support',
    'back compression support women',
    'back joint',
    'back pain compression',
    'back seam tights for women',
    'back support copper fit',
    'back support crossfit',
    'back support for standing',
    'back support for women workout',
    'back support garments',
    'back support xl',
    'back tights',
    'bad athletics',
    'bad boy tank top',
    'bad breath stomach',
    'bad support',
    'badass clothes for women',
    'badass gear',
    'badass lady',
    'ball between legs',
    'ball dresses for women evening 2017',
    'ball steep light',
    'ball total 90',
    'ballet gear',
    'ballet leg warmers',
    'ballet leggings for women',
    'ballet pants women',
    'ballet tight',
    'ballet tights',
    'ballet tights for men',
    'ballet tights for women',
    'ballet training equipment',
    'ballet underwear',
    'ballet underwear men',
    'band leggings',
    'barely bread brand',
    'base layer women capri',
    'baselayer pants women',
    'basjetball net',
    'baskeball games',
    'baskeball goal',
    'battle compression shorts',
    'battle sports shorts',
    'bball hoop',
    'beachbody leggings',
    'bed pants for women',
    'belly band surgery',
    'belly exercise',
    'belly fitness',
    'belly support leggings',
    'below knee leggings',
    'bermuda leggings',
    'bermuda leggings for women',
    'bermuda pants',
    'bermuda tights',
    'bermudas women',
    'best amazon leggings',
    'best basketball',
    'best black leggings for women',
    'best compression',
    'best compression leggings',
    'best compression leggings women',
    'best compression pants women',
    'best compression pantyhose high waist',
    'best compression shorts',
    'best compression tights',
    'best fitness',
    'best fitting leggings',
    'best generator',
    'best joint',
    'best joint hip and',
    'best joint support',
    'best karate',
    'best knee brace for skiing',
    'best knee support for running',
    'best leggings',
    'best leggings for running',
    'best leggings for winter',
    'best leggings for women',
    'best running leggings',
    'best running leggings for women',
    'best running pants',
    'best selling back brace'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
